$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '75.129.85'
$ws.Range('E2').Value = '  +2.03%  '

$ws.Range('D3').Value = '2.819.83'
$ws.Range('E3').Value = '  +7.90%  '

$ws.Range('E4').Value = '  +0.09%  '

$ws.Range('D5').Value = "'188.49"
$ws.Range('E5').Value = '  +2.63%  '

$ws.Range('D6').Value = "'595.86"
$ws.Range('E6').Value = '  +2.64%  '

$ws.Range('E7').Value = '  +0.09%  '

$ws.Range('E8').Value = '  +4.07%  '

$ws.Range('E9').Value = '  -1.37%  '

$ws.Range('D10').Value = '2.817.77'
$ws.Range('E10').Value = '  +7.77%  '

$ws.Range('E11').Value = '  -0.67%  '

$ws.Range('E12').Value = '  +3.17%  '

$ws.Range('D13').Value = "'4.85"

$ws.Range('D14').Value = '3.337.87'
$ws.Range('E14').Value = '  +8.50%  '

$ws.Range('D15').Value = '75.020.53'
$ws.Range('E15').Value = '  +2.05%  '

$ws.Range('D16').Value = "'0.0000188"
$ws.Range('E16').Value = '  +0.93%  '

$ws.Range('D17').Value = "'27.02"
$ws.Range('E17').Value = '  +4.33%  '

$ws.Range('D18').Value = '2.819.90'
$ws.Range('E18').Value = '  +8.04%  '

$ws.Range('D19').Value = "'8.95"
$ws.Range('E19').Value = '  -0.31%  '

$ws.Range('D20').Value = "'12.35"
$ws.Range('E20').Value = '  +4.52%  '

$ws.Range('D21').Value = "'376.76"
$ws.Range('E21').Value = '  +1.56%  '

$ws.Range('E22').Value = '  +0.84%  '

$ws.Range('D23').Value = "'4.09"
$ws.Range('E23').Value = '  +0.96%  '

$ws.Range('D24').Value = "'6.19"
$ws.Range('E24').Value = '  -0.20%  '

$ws.Range('D25').Value = "'0.999"
$ws.Range('E25').Value = '  -0.15%  '

$ws.Range('D26').Value = "'71.05"
$ws.Range('E26').Value = '  +2.24%  '

$ws.Range('B27').Value = 'WrappedeETH'
$ws.Range('C27').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D27').Value = '2.970.43'
$ws.Range('E27').Value = '  +9.47%  '

$ws.Range('B28').Value = 'NEARProtocol'
$ws.Range('C28').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D28').Value = "'4.18"
$ws.Range('E28').Value = '  +2.02%  '

$ws.Range('D29').Value = "'9.62"
$ws.Range('E29').Value = '  +4.51%  '

$ws.Range('E30').Value = '  +11.73%  '

$ws.Range('D31').Value = "'0.999"
$ws.Range('E31').Value = '  +0.09%  '

$ws.Range('D32').Value = "'515.97"
$ws.Range('E32').Value = '  +1.45%  '

$ws.Range('E33').Value = '  +1.46%  '

$ws.Range('D34').Value = "'7.81"
$ws.Range('E34').Value = '  +2.81%  '

$ws.Range('E35').Value = '  +4.11%  '

$ws.Range('E36').Value = '  -0.02%  '

$ws.Range('D37').Value = "'162.86"
$ws.Range('E37').Value = '  +1.26%  '

$ws.Range('D38').Value = "'20.06"
$ws.Range('E38').Value = '  +5.27%  '

$ws.Range('E39').Value = '  -0.94%  '

$ws.Range('D40').Value = "'19.43"
$ws.Range('E40').Value = '  +0.94%  '

$ws.Range('D41').Value = "'185.12"
$ws.Range('E41').Value = '  +18.02%  '

$ws.Range('E42').Value = '  +0.03%  '

$ws.Range('D43').Value = "'5.06"
$ws.Range('E43').Value = '  +4.66%  '

$ws.Range('E44').Value = '  +5.54%  '

$ws.Range('E45').Value = '  +1.16%  '

$ws.Range('E46').Value = '  +3.96%  '

$ws.Range('D47').Value = "'40.04"
$ws.Range('E47').Value = '  +4.06%  '

$ws.Range('D48').Value = "'2.33"
$ws.Range('E48').Value = '  +0.31%  '

$ws.Range('D49').Value = "'0.0859"
$ws.Range('E49').Value = '  -0.74%  '

$ws.Range('D50').Value = "'0.572"
$ws.Range('E50').Value = '  +9.16%  '

$ws.Range('D51').Value = "'3.73"
$ws.Range('E51').Value = '  +3.66%  '
